# Generate Report for Archive
#
# The localization status report moved from "Ready for handoff" to
# "In Translation" for the two e2e test files, on every sheet that
# surfaces a Status column (Overview's per-locale status columns, and
# each locale sheet's own Status column). The Status columns are also
# narrowed now that the new text is shorter.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) / de-de (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.42
$wsOverview.Columns.Item(6).ColumnWidth = 12.42

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.42

# --- de-de sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.42
